$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Convert E64:E66 bsecode values from text to true numbers
$ws.Range("E64").Value = 509930
$ws.Range("E65").Value = 590024
$ws.Range("E66").Value = 543220

# Append new rows 67-69 (next batch of screener results)
$ws.Range("A67").Value = "26/06/2024 07:45:44"
$ws.Range("B67").Value = 1
$ws.Range("C67").Value = "SUPREMEIND"
$ws.Range("D67").Value = "Supreme Industries Limited"
$ws.Range("E67").NumberFormat = "@"
$ws.Range("E67").Value = "509930"
$ws.Range("F67").Value = -1.58
$ws.Range("G67").Value = 5795
$ws.Range("H67").Value = 52862

$ws.Range("A68").Value = "26/06/2024 07:45:44"
$ws.Range("B68").Value = 2
$ws.Range("C68").Value = "FACT"
$ws.Range("D68").Value = "Fertilizers And Chemicals Travancore Limited"
$ws.Range("E68").NumberFormat = "@"
$ws.Range("E68").Value = "590024"
$ws.Range("F68").Value = 1.6
$ws.Range("G68").Value = 1015
$ws.Range("H68").Value = 2160639

$ws.Range("A69").Value = "26/06/2024 07:45:44"
$ws.Range("B69").Value = 3
$ws.Range("C69").Value = "MAXHEALTH"
$ws.Range("D69").Value = "Max Healthcare Institute Ltd"
$ws.Range("E69").NumberFormat = "@"
$ws.Range("E69").Value = "543220"
$ws.Range("F69").Value = -2.01
$ws.Range("G69").Value = 875.5
$ws.Range("H69").Value = 1270101
